$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$scratch.NumberFormat = "@"
$scratch.Value = "69.317.45"
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  -0.86%  "

$scratch.NumberFormat = "@"
$scratch.Value = "3.442.11"
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  -1.77%  "

$scratch.NumberFormat = "@"
$scratch.Value = "1.00"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.04%  "

$scratch.NumberFormat = "@"
$scratch.Value = "610.70"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +1.57%  "

$scratch.NumberFormat = "@"
$scratch.Value = "167.46"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -3.44%  "

$scratch.NumberFormat = "@"
$scratch.Value = "3.435.27"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -1.88%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.594"
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -2.40%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("E11").Value = "  -3.10%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.564"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -3.12%  "

$scratch.NumberFormat = "@"
$scratch.Value = "44.23"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -4.16%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.0000269"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -1.87%  "

$scratch.NumberFormat = "@"
$scratch.Value = "4.001.89"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -1.55%  "

$scratch.NumberFormat = "@"
$scratch.Value = "8.14"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -1.72%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$scratch.NumberFormat = "@"
$scratch.Value = "69.519.20"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$scratch.NumberFormat = "@"
$scratch.Value = "3.446.86"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$scratch.NumberFormat = "@"
$scratch.Value = "580.96"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -4.96%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.120"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.93%  "

$scratch.NumberFormat = "@"
$scratch.Value = "17.15"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.29%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.845"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -3.13%  "

$scratch.NumberFormat = "@"
$scratch.Value = "8.88"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -2.32%  "

$scratch.NumberFormat = "@"
$scratch.Value = "95.71"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.21%  "

$scratch.NumberFormat = "@"
$scratch.Value = "15.15"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -2.56%  "

$scratch.NumberFormat = "@"
$scratch.Value = "3.63"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -2.04%  "

$ws.Range("E27").Value = "  +0.04%  "

$scratch.NumberFormat = "@"
$scratch.Value = "2.42"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -5.56%  "

$scratch.NumberFormat = "@"
$scratch.Value = "32.72"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -4.44%  "

$scratch.NumberFormat = "@"
$scratch.Value = "8.61"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -3.87%  "

$scratch.NumberFormat = "@"
$scratch.Value = "7.80"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -3.68%  "

$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$scratch.NumberFormat = "@"
$scratch.Value = "1.24"
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -2.54%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$scratch.NumberFormat = "@"
$scratch.Value = "2.79"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -5.86%  "

$scratch.NumberFormat = "@"
$scratch.Value = "6.53"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -5.77%  "

$scratch.NumberFormat = "@"
$scratch.Value = "575.43"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -11.83%  "

$ws.Range("B36").Value = "Cosmos"
$ws.Range("C36").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$scratch.NumberFormat = "@"
$scratch.Value = "10.51"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -1.80%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$scratch.NumberFormat = "@"
$scratch.Value = "0.0476"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +0.19%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.0955"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -4.43%  "

$ws.Range("E39").Value = "  +0.39%  "

$scratch.NumberFormat = "@"
$scratch.Value = "55.89"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -0.92%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.140"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -1.77%  "

$scratch.NumberFormat = "@"
$scratch.Value = "3.14"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -12.40%  "

$scratch.NumberFormat = "@"
$scratch.Value = "3.236.44"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -2.61%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.0₃0683"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -0.87%  "

$scratch.NumberFormat = "@"
$scratch.Value = "31.05"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -3.91%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.294"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -5.13%  "

$scratch.NumberFormat = "@"
$scratch.Value = "2.76"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -5.58%  "

$scratch.NumberFormat = "@"
$scratch.Value = "2.40"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -6.14%  "

$ws.Range("E49").Value = "  -2.81%  "

$scratch.NumberFormat = "@"
$scratch.Value = "133.93"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -0.25%  "

$scratch.Clear()